$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: PointFive / VP EMEA -> Peter Reeve, 1st Interview, 2025-11-30 (45991)
$ws.Range("A5").Value = 677
$ws.Range("B5").Value = "PointFive"
$ws.Range("C5").Value = "VP EMEA"
$ws.Range("D5").Value = "Peter Reeve"
$ws.Range("E5").Value = "1st Interview"
$ws.Range("F5").Value = 45991

# Row 6: Cognition AI / Forward Deployed Engineer / Sales Engineer (UK) -> Simone Malekar, 1st Interview, 45991
$ws.Range("A6").Value = 702
$ws.Range("B6").Value = "Cognition AI"
$ws.Range("C6").Value = "Forward Deployed Engineer / Sales Engineer (UK)"
$ws.Range("D6").Value = "Simone Malekar"
$ws.Range("E6").Value = "1st Interview"
$ws.Range("F6").Value = 45991

# Row 7: Cognition AI / Forward Deployed Engineer / Sales Engineer (Middle East) -> Andrej Chomutovskij, CV Sent, 45983
$ws.Range("A7").Value = 714
$ws.Range("B7").Value = "Cognition AI"
$ws.Range("C7").Value = "Forward Deployed Engineer / Sales Engineer (Middle East)"
$ws.Range("D7").Value = "Andrej Chomutovskij"
$ws.Range("E7").Value = "CV Sent"
$ws.Range("F7").Value = 45983

# Row 8 (new): Cognition AI / Forward Deployed Engineer / Sales Engineer (Middle East) -> Oliver Waterman, 1st Interview, 45991
$ws.Range("A8").Value = 714
$ws.Range("B8").Value = "Cognition AI"
$ws.Range("C8").Value = "Forward Deployed Engineer / Sales Engineer (Middle East)"
$ws.Range("D8").Value = "Oliver Waterman"
$ws.Range("E8").Value = "1st Interview"
$ws.Range("F8").Value = 45991

# Row 9 (new): Cognition AI / Forward Deployed Engineer / Sales Engineer (Middle East) -> Aamer Mushtaq, CV Sent, 45983
$ws.Range("A9").Value = 714
$ws.Range("B9").Value = "Cognition AI"
$ws.Range("C9").Value = "Forward Deployed Engineer / Sales Engineer (Middle East)"
$ws.Range("D9").Value = "Aamer Mushtaq"
$ws.Range("E9").Value = "CV Sent"
$ws.Range("F9").Value = 45983

# Row 10 (new): Cognition AI / Forward Deployed Engineer / Sales Engineer (Middle East) -> JOSHUA TANNER, 1st Interview, 45991
$ws.Range("A10").Value = 714
$ws.Range("B10").Value = "Cognition AI"
$ws.Range("C10").Value = "Forward Deployed Engineer / Sales Engineer (Middle East)"
$ws.Range("D10").Value = "JOSHUA TANNER"
$ws.Range("E10").Value = "1st Interview"
$ws.Range("F10").Value = 45991

# Apply the date style (style index 2 in original workbook) to F8:F10 to match F5:F7 formatting
$ws.Range("F8:F10").NumberFormat = $ws.Range("F7").NumberFormat
